$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 23736.777
$ws.Range("I43").Value = 34233.0
$ws.Range("J43").Value = 18488.666
$ws.Range("K43").Value = 34233.0
$ws.Range("L43").Value = 18488.666
$ws.Range("M43").Value = -34164.0
$ws.Range("N43").Value = -18626.666
$ws.Range("H74").Value = 47137.207
$ws.Range("I74").Value = 57139.156
$ws.Range("J74").Value = 9129.8
$ws.Range("K74").Value = 57139.156
$ws.Range("L74").Value = 9129.8
$ws.Range("M74").Value = -56203.156
$ws.Range("N74").Value = -11001.8
$ws.Range("H77").Value = 47137.207
$ws.Range("I77").Value = 57139.156
$ws.Range("J77").Value = 9129.8
$ws.Range("K77").Value = 285695.78
$ws.Range("L77").Value = 45649.0
$ws.Range("M77").Value = -281015.78
$ws.Range("N77").Value = -55009.0
$ws.Range("H80").Value = 1058.7307
$ws.Range("I80").Value = 1073.5834
$ws.Range("J80").Value = 1046.0
$ws.Range("K80").Value = 3220.7502
$ws.Range("L80").Value = 3138.0
$ws.Range("M80").Value = -2222.7502
$ws.Range("N80").Value = -5134.0
$ws.Range("H83").Value = 1058.7307
$ws.Range("I83").Value = 1073.5834
$ws.Range("J83").Value = 1046.0
$ws.Range("K83").Value = 9662.2506
$ws.Range("L83").Value = 9414.0
$ws.Range("M83").Value = -4670.250599999999
$ws.Range("N83").Value = -19398.0
$ws.Range("H131").Value = 3482.875
$ws.Range("J131").Value = 23975.0
$ws.Range("L131").Value = 71925.0
$ws.Range("N131").Value = -82005.0
$ws.Range("H137").Value = 25000.0
$ws.Range("I137").Value = 0.0
$ws.Range("J137").Value = 25000.0
$ws.Range("K137").Value = 0.0
$ws.Range("L137").Value = ""
$ws.Range("M137").Value = ""
$ws.Range("N137").Value = -80100.0

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1550.4242
$ws.Range("I2").Value = 1610.52
$ws.Range("K2").Value = 1610.52
$ws.Range("M2").Value = -1497.52
$ws.Range("H22").Value = 2686.25
$ws.Range("I22").Value = 2686.25
$ws.Range("K22").Value = 2686.25
$ws.Range("M22").Value = -2387.25
$ws.Range("H32").Value = 2451.5454
$ws.Range("I32").Value = 2117.7114
$ws.Range("K32").Value = 2117.7114
$ws.Range("M32").Value = -1830.7114
$ws.Range("H97").Value = 998.44116
$ws.Range("I97").Value = 993.65625
$ws.Range("K97").Value = 993.65625
$ws.Range("M97").Value = -497.65625
$ws.Range("H116").Value = 1550.4242
$ws.Range("I116").Value = 1610.52
$ws.Range("K116").Value = 1610.52
$ws.Range("M116").Value = 683.48
$ws.Range("H124").Value = 20214.5
$ws.Range("J124").Value = 20214.5
$ws.Range("L124").Value = 20214.5
$ws.Range("N124").Value = -30034.5
$ws.Range("H125").Value = 28609.727
$ws.Range("J125").Value = 28609.727
$ws.Range("L125").Value = 28609.727
$ws.Range("N125").Value = -38449.727

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1550.4242
$ws.Range("I3").Value = 1610.52
$ws.Range("K3").Value = 1610.52
$ws.Range("M3").Value = -1496.52
$ws.Range("H75").Value = 32607.0
$ws.Range("I75").Value = 32607.0
$ws.Range("K75").Value = 32607.0
$ws.Range("M75").Value = -31671.0
$ws.Range("H78").Value = 32607.0
$ws.Range("I78").Value = 32607.0
$ws.Range("K78").Value = 97821.0
$ws.Range("M78").Value = -93141.0
$ws.Range("H130").Value = 62500.0
$ws.Range("J130").Value = 62500.0
$ws.Range("L130").Value = 62500.0
$ws.Range("N130").Value = -72540.0

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12890590.0
$ws.Range("I31").Value = 23898238.0
$ws.Range("J31").Value = 48332.668
$ws.Range("K31").Value = 23898238.0
$ws.Range("L31").Value = 48332.668
$ws.Range("M31").Value = -23897943.0
$ws.Range("N31").Value = -48922.668
$ws.Range("H34").Value = 12890590.0
$ws.Range("I34").Value = 23898238.0
$ws.Range("J34").Value = 48332.668
$ws.Range("K34").Value = 23898238.0
$ws.Range("L34").Value = 48332.668
$ws.Range("M34").Value = -23898036.0
$ws.Range("N34").Value = -48736.668
$ws.Range("H86").Value = 4470.0
$ws.Range("I86").Value = 4470.0
$ws.Range("K86").Value = 4470.0
$ws.Range("M86").Value = -3347.0
$ws.Range("H89").Value = 4470.0
$ws.Range("I89").Value = 4470.0
$ws.Range("K89").Value = 22350.0
$ws.Range("M89").Value = -16734.0
$ws.Range("H132").Value = 5255.5713
$ws.Range("I132").Value = 4964.8335
$ws.Range("K132").Value = 14894.5005
$ws.Range("M132").Value = -12364.5005

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 3849.5
$ws.Range("J39").Value = 4327.0
$ws.Range("L39").Value = 12981.0
$ws.Range("N39").Value = -13569.0
$ws.Range("H110").Value = 9994.75
$ws.Range("I110").Value = 9994.75
$ws.Range("K110").Value = 29984.25
$ws.Range("M110").Value = -25894.25
$ws.Range("H125").Value = 0.0
$ws.Range("J125").Value = 0.0
$ws.Range("L125").Value = ""
$ws.Range("N125").Value = ""
$ws.Range("H132").Value = 1344.3
$ws.Range("I132").Value = 1215.8889
$ws.Range("K132").Value = 10943.0001
$ws.Range("M132").Value = -8413.0001
$ws.Range("H140").Value = 11805.167
$ws.Range("I140").Value = 11805.167
$ws.Range("K140").Value = 35415.501
$ws.Range("M140").Value = -30235.501

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7537.0
$ws.Range("I70").Value = 6983.7646
$ws.Range("J70").Value = 8712.625
$ws.Range("K70").Value = 6983.7646
$ws.Range("L70").Value = 8712.625
$ws.Range("M70").Value = -6713.7646
$ws.Range("N70").Value = -9252.625
$ws.Range("H73").Value = 7537.0
$ws.Range("I73").Value = 6983.7646
$ws.Range("J73").Value = 8712.625
$ws.Range("K73").Value = 6983.7646
$ws.Range("L73").Value = 8712.625
$ws.Range("M73").Value = -6047.7646
$ws.Range("N73").Value = -10584.625
$ws.Range("H102").Value = 2699.5454
$ws.Range("I102").Value = 2699.5454
$ws.Range("K102").Value = 2699.5454
$ws.Range("M102").Value = -1077.5454
$ws.Range("H126").Value = 22234.25
$ws.Range("I126").Value = 28312.75
$ws.Range("K126").Value = 84938.25
$ws.Range("M126").Value = -82468.25
$ws.Range("H132").Value = 246675.97
$ws.Range("I132").Value = 259069.64
$ws.Range("K132").Value = 777208.92
$ws.Range("M132").Value = -774678.92
$ws.Range("H133").Value = 88500.0
$ws.Range("J133").Value = 88500.0
$ws.Range("L133").Value = 88500.0
$ws.Range("N133").Value = -98620.0

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3970.8076
$ws.Range("I132").Value = 3493.35
$ws.Range("J132").Value = 5562.3335
$ws.Range("K132").Value = 10480.05
$ws.Range("L132").Value = 16687.0005
$ws.Range("M132").Value = -7950.049999999999
$ws.Range("N132").Value = -21747.0005

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 25050.0
$ws.Range("I2").Value = 25050.0
$ws.Range("K2").Value = 25050.0
$ws.Range("M2").Value = -24938.0
$ws.Range("H96").Value = 2075.0
$ws.Range("I96").Value = 2075.0
$ws.Range("K96").Value = 2075.0
$ws.Range("M96").Value = -702.0
